$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D56: 13 -> 5
$ws.Range("D56").Value = 5

# A59: Memphis -> Florida
$ws.Range("A59").Value = "Florida"

# B62: Saint Johns -> Florida; D62: 5 -> 2
$ws.Range("B62").Value = "Florida"
$ws.Range("D62").Value = 2

# Row 64: remove the special border/fill style (now same as default, no style attrs)
# A64: Auburn -> Florida; C64: 5 -> 1; D64 stays 1 but loses its style
$ws.Range("A64").Value = "Florida"
$ws.Range("B64").Value = "Duke"
$ws.Range("C64").Value = 1
$ws.Range("D64").Value = 1
$ws.Range("A64:D64").Style = "Normal"

# Update the view: selection moves to A65 (below the last row of data)
$ws.Activate()
$ws.Range("A65").Select()
